# Reorder the EC (Estado de Cuenta) worker/period rows so each worker's
# two "Periodo Mora" rows (1805, then 1804) are grouped together, instead
# of being split into two separate blocks (all-1804 rows then all-1805 rows).
# Underlying worker/period/value data does not change, only its row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New row order (r16..r25): Tipo Doc | N Doc Trabajador | Nombre Trabajador | Periodo Mora | Valor Mora | Salario Basico
$data = @(
    @("CC", "1047367830", "ISABEL MARIA ACEVEDO MARQUEZ",    "1805", 31249, 781242),
    @("CC", "1047367830", "ISABEL MARIA ACEVEDO MARQUEZ",    "1804", 31249, 781242),
    @("CC", "73146457",   "VICENTE MENDOZA HERRERA",         "1805", 29509, 781242),
    @("CC", "73146457",   "VICENTE MENDOZA HERRERA",         "1804", 31249, 781242),
    @("CC", "1047425699", "DANIELA ALEJANDRA GARCIA CASTRO", "1805", 31249, 781242),
    @("CC", "1047425699", "DANIELA ALEJANDRA GARCIA CASTRO", "1804", 31249, 781242),
    @("CC", "45498230",   "MIRTA CECILIA GAVIRIA FRANCO",    "1805", 31249, 781242),
    @("CC", "45498230",   "MIRTA CECILIA GAVIRIA FRANCO",    "1804", 31249, 781242),
    @("CC", "45373331",   "CLARIBEL ALVAREZ GARCIA",         "1805", 31249, 781242),
    @("CC", "45373331",   "CLARIBEL ALVAREZ GARCIA",         "1804", 31249, 781242)
)

$row = 16
foreach ($line in $data) {
    $ws.Cells.Item($row, 2).Value = $line[0]   # B: Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $line[1]   # C: N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $line[2]   # D: Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $line[3]   # E: Periodo Mora
    $ws.Cells.Item($row, 6).Value = $line[4]   # F: Valor Mora
    $ws.Cells.Item($row, 7).Value = $line[5]   # G: Salario Basico
    $row = $row + 1
}

# Columns B:J use bestFit widths; the re-saved workbook widened them
# slightly (newer Excel build's font metrics). Re-apply the observed
# widths so column sizing matches the post-edit workbook.
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.333333333333332
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333336
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
